$d = $word.ActiveDocument

# -----------------------------------------------------------------
# 1) "Material Design Lite" paragraph: re-apply the same text via
#    Find/Replace so Word collapses the spell-checked "Lite" run
#    (and its surrounding proofErr markers) back into the
#    surrounding run, matching a clean re-save of that paragraph.
# -----------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "Material Design Lite to create a GUI in the form of a website.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Material Design Lite to create a GUI in the form of a website.", 2)

# -----------------------------------------------------------------
# 2) Remove the stray _GoBack bookmark that currently sits right
#    after "team's" in the "majority of the team's first coding
#    language" paragraph.
# -----------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

# -----------------------------------------------------------------
# 3) Append a brand-new closing paragraph to the end of the Summary
#    section, wrapped in parentheses, with a fresh _GoBack bookmark
#    placed where the original author's cursor last was (between the
#    hyphen and "work" of "framework").
# -----------------------------------------------------------------
$seg1 = "("
$seg2 = "Ultimately Python provided a reasonable middle ground between team members knowledge and the potential difficulty in setting up the relevant infrastructure and coding frame"
$seg3 = "-"
$seg4 = "work to carry out the project. This middle ground led to the logical choice in deciding to implement our project using python"
$seg5 = ")"
$fullText = $seg1 + $seg2 + $seg3 + $seg4 + $seg5

$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc - 1, $endOfDoc - 1)
$tail.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last.Range
$paraStart = $newPara.Start

$insertRange = $d.Range($paraStart, $paraStart)
$insertRange.InsertAfter($fullText)

# Split the freshly-typed text back into separate runs (matching the
# original diff) by dropping and removing zero-width bookmarks at
# each internal boundary -- this forces Word to keep the runs
# distinct even though they all share identical formatting.
$b1 = $paraStart + $seg1.Length
$b2 = $b1 + $seg2.Length
$b3 = $b2 + $seg3.Length
$b4 = $b3 + $seg4.Length

$d.Bookmarks.Add("_split1", $d.Range($b1, $b1))
$d.Bookmarks.Add("_split2", $d.Range($b2, $b2))
$d.Bookmarks.Add("_split3", $d.Range($b3, $b3))
$d.Bookmarks.Add("_split4", $d.Range($b4, $b4))

$d.Bookmarks("_split1").Delete()
$d.Bookmarks("_split2").Delete()
$d.Bookmarks("_split4").Delete()

# Keep the boundary between "-" and "work" alive as the real
# _GoBack bookmark (collapsed / zero-length, same as the source).
$d.Bookmarks("_split3").Delete()
$d.Bookmarks.Add("_GoBack", $d.Range($b3, $b3))

Write-Output "edit complete"
